# Apply the edit: change cell E1's text from "cm" to "cm_or_AO2"
# and update the active cell selection to E2 (matching the recorded sheet view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "cm_or_AO2"

$ws.Range("E2").Select()
